# Updated cryptos list on Sat Feb 17 16:40:20 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "349.54") are briefly forced to Text format so they are stored as
# strings (matching the source data, which is all inlineStr), then the style
# is reset back to Normal so no stray number-format/style is left behind.

$ws.Range("D2").Value = "50.813.84"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "2.748.46"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "3.176.41"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "2.752.09"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "50.753.51"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -3.30%  "
$ws.Range("E28").Value = "  +12.71%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "51.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "2.077.50"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.909"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.92%  "
